$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 522.1818
$ws.Range("I19").Value = 481.2
$ws.Range("J19").Value = 556.3333
$ws.Range("K19").Value = 481.2
$ws.Range("L19").Value = 556.3333
$ws.Range("M19").Value = -306.2
$ws.Range("N19").Value = -906.3333

$ws.Range("H29").Value = 997
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws.Range("H41").Value = 914.7
$ws.Range("I41").Value = 736.5
$ws.Range("K41").Value = 736.5
$ws.Range("M41").Value = -296.5

$ws.Range("H43").Value = 9767.1
$ws.Range("I43").Value = 997
$ws.Range("J43").Value = 11959.625
$ws.Range("K43").Value = 997
$ws.Range("L43").Value = 11959.625
$ws.Range("M43").Value = -928
$ws.Range("N43").Value = -12097.625

$ws.Range("H58").Value = 2884.8667
$ws.Range("I58").Value = 616.25
$ws.Range("J58").Value = 3709.818
$ws.Range("K58").Value = 1848.75
$ws.Range("L58").Value = 11129.454
$ws.Range("M58").Value = -1698.75
$ws.Range("N58").Value = -11429.454

$ws.Range("H69").Value = 10000
$ws.Range("I69").Value = 10000
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 30000
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -29126
$ws.Range("N69").ClearContents()

$ws.Range("H70").Value = 15385983
$ws.Range("I70").Value = 50000724
$ws.Range("J70").Value = 1652.7778
$ws.Range("K70").Value = 150002172
$ws.Range("L70").Value = 4958.3334
$ws.Range("M70").Value = -150001902
$ws.Range("N70").Value = -5498.3334

$ws.Range("H72").Value = 10000
$ws.Range("I72").Value = 10000
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 90000
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -85632
$ws.Range("N72").ClearContents()

$ws.Range("H73").Value = 15385983
$ws.Range("I73").Value = 50000724
$ws.Range("J73").Value = 1652.7778
$ws.Range("K73").Value = 150002172
$ws.Range("L73").Value = 4958.3334
$ws.Range("M73").Value = -150001236
$ws.Range("N73").Value = -6830.3334

$ws.Range("H106").Value = 14870.833
$ws.Range("I106").Value = 14845
$ws.Range("J106").Value = 15000
$ws.Range("K106").Value = 14845
$ws.Range("L106").Value = 15000
$ws.Range("M106").Value = -14214
$ws.Range("N106").Value = -16262

$ws.Range("H132").Value = 2086613.8
$ws.Range("J132").Value = 14288978
$ws.Range("L132").Value = 42866934
$ws.Range("N132").Value = -42871994

$ws.Range("H137").Value = 10221.04
$ws.Range("I137").Value = 18949.834
$ws.Range("J137").Value = 2163.6924
$ws.Range("K137").Value = 56849.50199999999
$ws.Range("L137").Value = 6491.0772
$ws.Range("M137").Value = -54299.50199999999
$ws.Range("N137").Value = -11591.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 95.25
$ws.Range("I4").Value = 95.25
$ws.Range("K4").Value = 95.25
$ws.Range("M4").Value = 20.75

$ws.Range("H122").Value = 972609.5
$ws.Range("J122").Value = 2732389.5
$ws.Range("L122").Value = 8197168.5
$ws.Range("N122").Value = -8202068.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws.Range("H31").Value = 10844.823
$ws.Range("I31").Value = 12866.083
$ws.Range("J31").Value = 5993.8
$ws.Range("K31").Value = 12866.083
$ws.Range("L31").Value = 5993.8
$ws.Range("M31").Value = -12571.083
$ws.Range("N31").Value = -6583.8

$ws.Range("H34").Value = 10844.823
$ws.Range("I34").Value = 12866.083
$ws.Range("J34").Value = 5993.8
$ws.Range("K34").Value = 12866.083
$ws.Range("L34").Value = 5993.8
$ws.Range("M34").Value = -12664.083
$ws.Range("N34").Value = -6397.8

$ws.Range("H41").Value = 16541.666
$ws.Range("I41").Value = 312.5
$ws.Range("K41").Value = 312.5
$ws.Range("M41").Value = 115.5

$ws.Range("H50").Value = 68998.5
$ws.Range("J50").Value = 68998.5
$ws.Range("L50").Value = 68998.5
$ws.Range("N50").Value = -70248.5

$ws.Range("H51").Value = 24997.25
$ws.Range("I51").Value = 10000
$ws.Range("J51").Value = 29996.334
$ws.Range("K51").Value = 10000
$ws.Range("L51").Value = 29996.334
$ws.Range("M51").Value = -9264
$ws.Range("N51").Value = -31468.334

$ws.Range("H59").Value = 72000
$ws.Range("J59").Value = 72000
$ws.Range("L59").Value = 72000
$ws.Range("N59").Value = -74290

$ws.Range("H60").Value = 42305.5
$ws.Range("J60").Value = 42305.5
$ws.Range("L60").Value = 42305.5
$ws.Range("N60").Value = -43327.5

$ws.Range("H61").Value = 24997.25
$ws.Range("I61").Value = 10000
$ws.Range("J61").Value = 29996.334
$ws.Range("K61").Value = 10000
$ws.Range("L61").Value = 29996.334
$ws.Range("M61").Value = -9652
$ws.Range("N61").Value = -30692.334

$ws.Range("H107").Value = 3369.6829
$ws.Range("I107").Value = 4301.5806
$ws.Range("J107").Value = 480.8
$ws.Range("K107").Value = 4301.5806
$ws.Range("L107").Value = 480.8
$ws.Range("M107").Value = -2381.5806
$ws.Range("N107").Value = -4320.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 550.5
$ws.Range("J23").Value = 529
$ws.Range("L23").Value = 1587
$ws.Range("N23").Value = -2057

$ws.Range("H34").Value = 2009.1818
$ws.Range("I34").Value = 1337.625
$ws.Range("J34").Value = 3800
$ws.Range("K34").Value = 4012.875
$ws.Range("L34").Value = 11400
$ws.Range("M34").Value = -3928.875
$ws.Range("N34").Value = -11568

$ws.Range("H39").Value = 300
$ws.Range("I39").Value = 300
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 900
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -606
$ws.Range("N39").ClearContents()

$ws.Range("H55").Value = 11055.566
$ws.Range("J55").Value = 11630.821
$ws.Range("L55").Value = 34892.463
$ws.Range("N55").Value = -35246.463

$ws.Range("H122").Value = 4190.6665
$ws.Range("J122").Value = 4553.421
$ws.Range("L122").Value = 40980.789
$ws.Range("N122").Value = -45880.789

$ws.Range("H131").Value = 1487.404
$ws.Range("I131").Value = 855
$ws.Range("J131").Value = 1521.0426
$ws.Range("K131").Value = 2565
$ws.Range("L131").Value = 4563.1278
$ws.Range("M131").Value = 2475
$ws.Range("N131").Value = -14643.1278

$ws.Range("H133").Value = 12497.5
$ws.Range("I133").Value = 4997
$ws.Range("K133").Value = 14991
$ws.Range("M133").Value = -9931

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 7236.8696
$ws.Range("I102").Value = 8488
$ws.Range("K102").Value = 8488
$ws.Range("M102").Value = -6866

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2508.5625
$ws.Range("I82").Value = 2789.5
$ws.Range("J82").Value = 2040.3334
$ws.Range("K82").Value = 2789.5
$ws.Range("L82").Value = 2040.3334
$ws.Range("M82").Value = -2428.5
$ws.Range("N82").Value = -2762.3334

$ws.Range("H85").Value = 2508.5625
$ws.Range("I85").Value = 2789.5
$ws.Range("J85").Value = 2040.3334
$ws.Range("K85").Value = 2789.5
$ws.Range("L85").Value = 2040.3334
$ws.Range("M85").Value = -1541.5
$ws.Range("N85").Value = -4536.3334

$ws.Range("H122").Value = 4701.4414
$ws.Range("I122").Value = 4526.2
$ws.Range("K122").Value = 13578.6
$ws.Range("M122").Value = -11128.6

$ws.Range("H132").Value = 879526.75
$ws.Range("I132").Value = 1356261.5
$ws.Range("J132").Value = 5513
$ws.Range("K132").Value = 4068784.5
$ws.Range("L132").Value = 16539
$ws.Range("M132").Value = -4066254.5
$ws.Range("N132").Value = -21599

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 383562.12
$ws.Range("I62").Value = 571105.7
$ws.Range("J62").Value = 8475
$ws.Range("K62").Value = 571105.7
$ws.Range("L62").Value = 8475
$ws.Range("M62").Value = -570481.7
$ws.Range("N62").Value = -9723

$ws.Range("H65").Value = 383562.12
$ws.Range("I65").Value = 571105.7
$ws.Range("J65").Value = 8475
$ws.Range("K65").Value = 2855528.5
$ws.Range("L65").Value = 42375
$ws.Range("M65").Value = -2852408.5
$ws.Range("N65").Value = -48615

$ws.Range("H95").Value = 172010000
$ws.Range("J95").Value = 172010000
$ws.Range("L95").Value = 172010000
$ws.Range("N95").Value = -172015492
